# Roll the quarterly reporting window forward by one quarter:
# drop the oldest quarter column (فصل دوم منتهی به 1399/06) and append the
# newest quarter column (فصل چهارم منتهی به 1401/12), shifting every
# existing quarter's data one column to the left (E<-F, F<-G, ..., M<-N)
# and filling the freed rightmost column (N) with the new quarter's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Shift-RowLeft {
    param($Row, $NewValue)

    $cols = @("E", "F", "G", "H", "I", "J", "K", "L", "M", "N")

    # Read the existing 10 values (columns E..N) left to right first,
    # since writes below would otherwise clobber values still to be read.
    $values = @()
    foreach ($col in $cols) {
        $values += , ($ws.Range("$col$Row").Value())
    }

    # Shift everything one column to the left ...
    for ($i = 0; $i -lt 9; $i++) {
        $ws.Range("$($cols[$i])$Row").Value = $values[$i + 1]
    }

    # ... and place the newest quarter's value in the freed last column.
    $ws.Range("N$Row").Value = $NewValue
}

# Header rows: quarter labels
Shift-RowLeft 8 "فصل چهارم منتهی به 1401/12"
Shift-RowLeft 24 "فصل چهارم منتهی به 1401/12"

# Expense rows (each total recomputed as the sum of the new quarter's
# category values, matching the existing "جمع" row semantics)
Shift-RowLeft 10 857
Shift-RowLeft 11 0
Shift-RowLeft 12 0
Shift-RowLeft 13 119528
Shift-RowLeft 14 0
Shift-RowLeft 15 0
Shift-RowLeft 16 1494
Shift-RowLeft 17 236977
Shift-RowLeft 18 0
Shift-RowLeft 19 70423

# "جمع" (Total) row
Shift-RowLeft 20 429279

# Personnel counts
Shift-RowLeft 26 574
Shift-RowLeft 27 656
